# Fix the broken domain used to build the PDF links in column E
# (http://asiaameircalatina.org.ar -> http://www.asiaamericalatina.org)
# and update the active selection as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old (misspelled) host across all formulas/values on the sheet.
# Using Cells.Replace keeps the existing shared-formula grouping intact,
# it just substitutes the literal text inside each formula.
$null = $ws.Cells.Replace("http://asiaameircalatina.org.ar/docs/", "http://www.asiaamericalatina.org/docs/")

# Restore the selection recorded in the saved file.
$ws.Range("F5").Select()
